# Training Log.xlsx - "project proposal 2nd submission"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Structural edits -------------------------------------------------
# Insert a new column B ("lookahead minutes") shifting dataset..notes right.
$ws.Columns("B:B").Insert()

# Insert two new rows at 5:6 (pushes the old row 5 "with inverse" down to row 7),
# giving us rows 5 and 6 to fill with new data.
$ws.Rows("5:6").Insert()

# That also pushed the old rows 13/14 down to 15/16; remove two now-empty rows
# from the gap above them so they land back on rows 13/14.
$ws.Rows("8:9").Delete()

# --- Row 5 (new) -----------------------------------------------------
$ws.Range("C5").Value2 = "single sided"
$ws.Range("D5").Value2 = "MSE"
$ws.Range("E5").Value2 = 20
$ws.Range("F5").Value2 = 0.1579
$ws.Range("G5").Value2 = 0.1567
$ws.Range("H5").Value2 = 0.256
$ws.Range("I5").Value2 = 0.22
$ws.Range("J5").Value2 = "repeat to check consistency"

# --- Row 6 (was row 5 pre-insert: "with inverse"/"MAE"/10) ---------------
$ws.Range("C6").Value2 = "with inverse"
$ws.Range("D6").Value2 = "MSE"
$ws.Range("E6").Value2 = 10

# --- Row 7 (new) -----------------------------------------------------------
$ws.Range("C7").Value2 = "with inverse"
$ws.Range("D7").Value2 = "MSE"
$ws.Range("E7").Value2 = 20

# --- Header row (row 1) -------------------------------------------------
$ws.Range("B1").Value2 = "lookahead minutes"

# --- Row 2 ---------------------------------------------------------------
$ws.Range("B2").Value2 = 10

# --- Old rows 13/14 (TODO:/X notes) removed, keep two styled blank cells -
$ws.Range("A13:J14").ClearContents()
$ws.Rows("14:14").Delete()
$ws.Range("B13").Font.Bold = $true

# --- Column widths (best effort) -----------------------------------------
# The saved XML "width" = ColumnWidth(chars) + 5/MDW (a fixed ~0.8333 padding
# this engine always adds once a custom width is set, MDW=6px here). We
# subtract that padding back out so the serialized width matches the target.
$pad = 5.0/6.0
$ws.Columns("A:A").ColumnWidth = 24.83203125 - $pad
$ws.Columns("B:B").ColumnWidth = 10 - $pad
$ws.Columns("C:C").ColumnWidth = 9.83203125 - $pad
$ws.Columns("D:D").ColumnWidth = 11.5 - $pad
$ws.Columns("E:E").ColumnWidth = 6.1640625 - $pad
$ws.Columns("F:F").ColumnWidth = 7.83203125 - $pad
$ws.Columns("G:G").ColumnWidth = 10.5 - $pad
$ws.Columns("H:H").ColumnWidth = 12.1640625 - $pad
$ws.Columns("I:I").ColumnWidth = 7.1640625 - $pad
$ws.Columns("J:J").ColumnWidth = 32 - $pad

# --- View / selection state ----------------------------------------------
$excel.ActiveWindow.Zoom = 150
$ws.Range("E7").Select()
